$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 545-653: Fila, Fecha(D), Calidad(L), Volumen(M), PrecioMin(N), PrecioMax(O), PrecioProm(P), PrecioKg(S)
$data = @(
    ,@(545, 44637, "Especial", 400, 12000, 12500, 12250, 1750)
    ,@(546, 44637, "Primera", 500, 10000, 10500, 10250, 1464)
    ,@(547, 44637, "Segunda", 400, 8000, 8500, 8250, 1179)
    ,@(548, 44208, "Especial", 300, 15500, 16000, 15750, 2250)
    ,@(549, 44208, "Primera", 260, 13500, 14000, 13750, 1964)
    ,@(550, 44445, "Primera", 100, 26000, 27000, 26500, 3786)
    ,@(551, 44445, "Segunda", 100, 20000, 21000, 20500, 2929)
    ,@(552, 44355, "Especial", 400, 17500, 18000, 17750, 2536)
    ,@(553, 44355, "Primera", 340, 15500, 16000, 15750, 2250)
    ,@(554, 44355, "Segunda", 260, 11500, 12000, 11750, 1679)
    ,@(555, 44530, "Especial", 400, 12500, 13000, 12750, 1821)
    ,@(556, 44530, "Primera", 360, 10500, 11000, 10750, 1536)
    ,@(557, 44530, "Segunda", 280, 8500, 9000, 8750, 1250)
    ,@(558, 44483, "Especial", 400, 14000, 14500, 14250, 2036)
    ,@(559, 44483, "Primera", 500, 12000, 12500, 12250, 1750)
    ,@(560, 44483, "Segunda", 300, 9000, 9500, 9250, 1321)
    ,@(561, 44294, "Especial", 300, 14500, 15000, 14750, 2107)
    ,@(562, 44294, "Primera", 300, 12500, 13000, 12750, 1821)
    ,@(563, 44294, "Segunda", 240, 10500, 11000, 10750, 1536)
    ,@(564, 44617, "Especial", 400, 12000, 12500, 12250, 1750)
    ,@(565, 44617, "Primera", 500, 10000, 10500, 10250, 1464)
    ,@(566, 44617, "Segunda", 400, 8000, 8500, 8250, 1179)
    ,@(567, 44489, "Especial", 300, 13500, 14000, 13750, 1964)
    ,@(568, 44489, "Primera", 500, 11500, 12000, 11750, 1679)
    ,@(569, 44489, "Segunda", 360, 9500, 10000, 9750, 1393)
    ,@(570, 44264, "Especial", 500, 15500, 16000, 15750, 2250)
    ,@(571, 44264, "Primera", 360, 13500, 14000, 13750, 1964)
    ,@(572, 44264, "Segunda", 200, 10500, 11000, 10750, 1536)
    ,@(573, 44232, "Especial", 240, 16000, 16500, 16250, 2321)
    ,@(574, 44232, "Primera", 240, 14000, 14500, 14250, 2036)
    ,@(575, 44232, "Segunda", 240, 12000, 12500, 12250, 1750)
    ,@(576, 44330, "Especial", 240, 19500, 20000, 19750, 2821)
    ,@(577, 44330, "Primera", 300, 16500, 17000, 16750, 2393)
    ,@(578, 44330, "Segunda", 240, 12000, 12500, 12250, 1750)
    ,@(579, 44504, "Especial", 400, 12000, 12500, 12250, 1750)
    ,@(580, 44504, "Primera", 500, 10000, 10500, 10250, 1464)
    ,@(581, 44504, "Segunda", 400, 8000, 8500, 8250, 1179)
    ,@(582, 44257, "Especial", 400, 15500, 16000, 15750, 2250)
    ,@(583, 44257, "Primera", 300, 13500, 14000, 13750, 1964)
    ,@(584, 44257, "Segunda", 200, 10500, 11000, 10750, 1536)
    ,@(585, 44301, "Especial", 240, 14500, 15000, 14750, 2107)
    ,@(586, 44301, "Primera", 300, 12500, 13000, 12750, 1821)
    ,@(587, 44301, "Segunda", 240, 10500, 11000, 10750, 1536)
    ,@(588, 44487, "Especial", 300, 13500, 14000, 13750, 1964)
    ,@(589, 44487, "Primera", 400, 11500, 12000, 11750, 1679)
    ,@(590, 44487, "Segunda", 300, 9500, 10000, 9750, 1393)
    ,@(591, 44174, "Especial", 240, 17500, 18000, 17750, 2536)
    ,@(592, 44174, "Primera", 240, 15500, 16000, 15750, 2250)
    ,@(593, 44174, "Segunda", 200, 13500, 14000, 13750, 1964)
    ,@(594, 44200, "Especial", 240, 15500, 16000, 15750, 2250)
    ,@(595, 44200, "Primera", 300, 13500, 14000, 13750, 1964)
    ,@(596, 44200, "Segunda", 260, 11500, 12000, 11750, 1679)
    ,@(597, 44385, "Especial", 240, 24500, 25000, 24750, 3536)
    ,@(598, 44385, "Primera", 160, 22500, 23000, 22750, 3250)
    ,@(599, 44385, "Segunda", 200, 16500, 17000, 16750, 2393)
    ,@(600, 44236, "Especial", 300, 16500, 17000, 16750, 2393)
    ,@(601, 44236, "Primera", 240, 13500, 14000, 13750, 1964)
    ,@(602, 44236, "Segunda", 200, 11500, 12000, 11750, 1679)
    ,@(603, 44221, "Primera", 300, 14500, 15000, 14750, 2107)
    ,@(604, 44221, "Segunda", 240, 11500, 12000, 11750, 1679)
    ,@(605, 44413, "Especial", 160, 26500, 27000, 26750, 3821)
    ,@(606, 44413, "Primera", 200, 21500, 22000, 21750, 3107)
    ,@(607, 44413, "Segunda", 200, 17500, 18000, 17750, 2536)
    ,@(608, 44229, "Especial", 300, 16500, 17000, 16750, 2393)
    ,@(609, 44229, "Primera", 280, 14500, 15000, 14750, 2107)
    ,@(610, 44229, "Segunda", 240, 12500, 13000, 12750, 1821)
    ,@(611, 44214, "Especial", 240, 17500, 18000, 17750, 2536)
    ,@(612, 44214, "Primera", 240, 15500, 16000, 15750, 2250)
    ,@(613, 44214, "Segunda", 240, 12500, 13000, 12750, 1821)
    ,@(614, 44610, "Especial", 400, 11500, 12000, 11750, 1679)
    ,@(615, 44610, "Primera", 400, 9500, 10000, 9750, 1393)
    ,@(616, 44610, "Segunda", 300, 7500, 8000, 7750, 1107)
    ,@(617, 44312, "Especial", 240, 15500, 16000, 15750, 2250)
    ,@(618, 44312, "Primera", 300, 13500, 14000, 13750, 1964)
    ,@(619, 44312, "Segunda", 300, 10500, 11000, 10750, 1536)
    ,@(620, 44399, "Especial", 240, 24500, 25000, 24750, 3536)
    ,@(621, 44399, "Primera", 300, 19500, 20000, 19750, 2821)
    ,@(622, 44399, "Segunda", 200, 16500, 17000, 16750, 2393)
    ,@(623, 44522, "Especial", 360, 12500, 13000, 12750, 1821)
    ,@(624, 44522, "Primera", 300, 10500, 11000, 10750, 1536)
    ,@(625, 44522, "Segunda", 200, 8500, 9000, 8750, 1250)
    ,@(626, 44543, "Especial", 300, 11500, 12000, 11750, 1679)
    ,@(627, 44543, "Primera", 300, 9500, 10000, 9750, 1393)
    ,@(628, 44543, "Segunda", 300, 7500, 8000, 7750, 1107)
    ,@(629, 44167, "Especial", 500, 15000, 16000, 15500, 2214)
    ,@(630, 44277, "Especial", 240, 13500, 14000, 13750, 1964)
    ,@(631, 44277, "Primera", 300, 11500, 12000, 11750, 1679)
    ,@(632, 44277, "Segunda", 240, 9500, 10000, 9750, 1393)
    ,@(633, 44390, "Especial", 400, 22000, 23000, 22500, 3214)
    ,@(634, 44390, "Primera", 300, 20000, 21000, 20500, 2929)
    ,@(635, 44390, "Segunda", 240, 18000, 19000, 18500, 2643)
    ,@(636, 44498, "Especial", 300, 12500, 13000, 12750, 1821)
    ,@(637, 44498, "Primera", 300, 10500, 11000, 10750, 1536)
    ,@(638, 44498, "Segunda", 300, 8500, 9000, 8750, 1250)
    ,@(639, 44179, "Especial", 240, 15500, 16000, 15750, 2250)
    ,@(640, 44179, "Primera", 240, 13500, 14000, 13750, 1964)
    ,@(641, 44179, "Segunda", 240, 11500, 12000, 11750, 1679)
    ,@(642, 44595, "Especial", 400, 11500, 12000, 11750, 1679)
    ,@(643, 44595, "Primera", 400, 9500, 10000, 9750, 1393)
    ,@(644, 44595, "Segunda", 360, 7500, 8000, 7750, 1107)
    ,@(645, 44628, "Especial", 520, 11000, 12000, 11500, 1643)
    ,@(646, 44628, "Primera", 400, 9000, 10000, 9500, 1357)
    ,@(647, 44628, "Segunda", 300, 7000, 8000, 7500, 1071)
    ,@(648, 44544, "Especial", 400, 11000, 12000, 11500, 1643)
    ,@(649, 44544, "Primera", 300, 9000, 10000, 9500, 1357)
    ,@(650, 44544, "Segunda", 240, 7500, 8000, 7750, 1107)
    ,@(651, 44160, "Especial", 400, 14500, 15000, 14750, 2107)
    ,@(652, 44160, "Primera", 300, 12500, 13000, 12750, 1821)
    ,@(653, 44160, "Segunda", 240, 10500, 11000, 10750, 1536)
)

# Static / boilerplate column values shared by every data row in this block
# NOTE: use Value2 (not Value) to read - this runtime's .Value getter does not
# resolve to the underlying scalar when read back into a variable.
$colA = $ws.Range("A650").Value2
$colB = $ws.Range("B650").Value2
$colC = $ws.Range("C650").Value2
$colE = $ws.Range("E650").Value2
$colF = $ws.Range("F650").Value2
$colG = $ws.Range("G650").Value2
$colH = $ws.Range("H650").Value2
$colI = $ws.Range("I650").Value2
$colJ = $ws.Range("J650").Value2
$colK = $ws.Range("K650").Value2
$colQ = $ws.Range("Q650").Value2
$colR = $ws.Range("R650").Value2
$colT = $ws.Range("T650").Value2
$dateFormat = $ws.Range("D650").NumberFormat

foreach ($row in $data) {
    $r = $row[0]
    if ($r -gt 650) {
        # Brand-new row: populate every column, not just the variable ones
        $ws.Range("A$r").Value = $colA
        $ws.Range("B$r").Value = $colB
        $ws.Range("C$r").Value = $colC
        $ws.Range("D$r").NumberFormat = $dateFormat
        $ws.Range("E$r").Value = $colE
        $ws.Range("F$r").Value = $colF
        $ws.Range("G$r").Value = $colG
        $ws.Range("H$r").Value = $colH
        $ws.Range("I$r").Value = $colI
        $ws.Range("J$r").Value = $colJ
        $ws.Range("K$r").Value = $colK
        $ws.Range("Q$r").Value = $colQ
        $ws.Range("R$r").Value = $colR
        $ws.Range("T$r").Value = $colT
    }
    $ws.Range("D$r").Value = $row[1]
    $ws.Range("L$r").Value = $row[2]
    $ws.Range("M$r").Value = $row[3]
    $ws.Range("N$r").Value = $row[4]
    $ws.Range("O$r").Value = $row[5]
    $ws.Range("P$r").Value = $row[6]
    $ws.Range("S$r").Value = $row[7]
}

Write-Host "Done: updated rows 545-650 and appended rows 651-653"
